$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the formatting of the existing header row
# (copy G1's style/formatting into H1, then set the new header text)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the data values for the new "Save" column (plain numbers, default style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
